# Updated cryptos list values per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.953.02'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '1.637.85'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  +0.98%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0636'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.63'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0795'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = '1.864.45'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = '1.637.45'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  -1.78%  '
$ws.Range("D16").Value = '0.0₃0756'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '25.958.55'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.87'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.25'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("E27").Value = '  +2.65%  '
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.54'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.20%  '
$ws.Range("E35").Value = '  +1.65%  '
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").Value = '1.138.96'
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.546'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.42'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("D44").Value = '1.774.18'
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '0.0₆0115'
$ws.Range("E45").Value = '  +8.97%  '
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("E47").Value = '  +2.89%  '
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.415'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.65'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0961'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.32%  '
